$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Row 5 (tester1@gmail.com / tester123) moves into row 6, overwriting what was
# there (demo4@example.com / test1234); row 5 becomes empty.
$ws.Range("A6").Value = "tester1@gmail.com"
$ws.Range("B6").Value = "tester123"
$ws.Range("A5:B5").ClearContents()

# A new row 9 is appended, duplicating the last row's content
# (demo4@example.com / test1234).
$ws.Range("A9").Value = "demo4@example.com"
$ws.Range("B9").Value = "test1234"

# Update the selection to match the recorded cursor position after the edit.
$null = $ws.Range("D10").Select()
